$wb = $excel.ActiveWorkbook

# --- Rename header cells on existing sheets -----------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ------
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast.Name = "PO Forecast"

# Copy the header formatting (bold, centered, bordered) from an existing
# header cell so the new sheet reuses the same cell style.
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy the date-number-format cell style down column A for the data rows.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A12").PasteSpecial(-4122)

$wsForecast.Range("A2").Value = 45431.99999999999
$wsForecast.Range("B2").Value = 22
$wsForecast.Range("C2").Value = 17.78008568298561
$wsForecast.Range("D2").Value = 25.28507986642154

$wsForecast.Range("A3").Value = 45480.99999999999
$wsForecast.Range("B3").Value = 12
$wsForecast.Range("C3").Value = 8.496288044944944
$wsForecast.Range("D3").Value = 16.24216869171887

$wsForecast.Range("A4").Value = 45543.99999999999
$wsForecast.Range("B4").Value = 0
$wsForecast.Range("C4").Value = -3.827895561742556
$wsForecast.Range("D4").Value = 4.075095712190635

$wsForecast.Range("A5").Value = 45550.99999999999
$wsForecast.Range("B5").Value = 0
$wsForecast.Range("C5").Value = -4.969103366917676
$wsForecast.Range("D5").Value = 2.446253781626911

$wsForecast.Range("A6").Value = 45557.99999999999
$wsForecast.Range("B6").Value = 0
$wsForecast.Range("C6").Value = -6.398747354236812
$wsForecast.Range("D6").Value = 1.112956190327953

$wsForecast.Range("A7").Value = 45564.99999999999
$wsForecast.Range("B7").Value = 0
$wsForecast.Range("C7").Value = -7.502813371030983
$wsForecast.Range("D7").Value = -0.1452932110955633

$wsForecast.Range("A8").Value = 45571.99999999999
$wsForecast.Range("B8").Value = 0
$wsForecast.Range("C8").Value = -8.976805686700454
$wsForecast.Range("D8").Value = -1.363842014881746

$wsForecast.Range("A9").Value = 45578.99999999999
$wsForecast.Range("B9").Value = 0
$wsForecast.Range("C9").Value = -10.26756134579615
$wsForecast.Range("D9").Value = -2.66563346326357

$wsForecast.Range("A10").Value = 45585.99999999999
$wsForecast.Range("B10").Value = 0
$wsForecast.Range("C10").Value = -11.74558316006116
$wsForecast.Range("D10").Value = -4.398807760805224

$wsForecast.Range("A11").Value = 45592.99999999999
$wsForecast.Range("B11").Value = 0
$wsForecast.Range("C11").Value = -13.20908509505201
$wsForecast.Range("D11").Value = -5.443654497900702

$wsForecast.Range("A12").Value = 45599.99999999999
$wsForecast.Range("B12").Value = 0
$wsForecast.Range("C12").Value = -14.55878994278505
$wsForecast.Range("D12").Value = -6.644946254111168


# Restore the original active sheet (adding a sheet makes it active).
$wsWeekly.Activate()
